$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1034
$ws.Range("E2").Value = 73
$ws.Range("F2").Value = 73
$ws.Range("G2").Value = 73
$ws.Range("H2").Value = 52
$ws.Range("I2").Value = 52
$ws.Range("K2").Value = 1010
$ws.Range("L2").Value = 355
$ws.Range("M2").Value = 655
$ws.Range("N2").Value = 655
$ws.Range("P2").Value = 128
$ws.Range("Q2").Value = 55
$ws.Range("R2").Value = -68
$ws.Range("S2").Value = -23
$ws.Range("T2").Value = 93
$ws.Range("U2").Value = -38
$ws.Range("V2").Value = 133
$ws.Range("W2").Value = 7.1
$ws.Range("X2").Value = 5.06
$ws.Range("Y2").Value = 8.31
$ws.Range("Z2").Value = 5.34
$ws.Range("AA2").Value = 54.12
$ws.Range("AB2").Value = 333.34
$ws.Range("AC2").Value = 201
$ws.Range("AD2").Value = 10.91
$ws.Range("AE2").Value = 2518
$ws.Range("AF2").Value = 0.87
$ws.Range("AG2").Value = 39
$ws.Range("AH2").Value = 1.77
$ws.Range("AI2").Value = 19.34
$ws.Range("AJ2").Value = 26027444

# Row 3
$ws.Range("D3").Value = 948
$ws.Range("E3").Value = 48
$ws.Range("F3").Value = 48
$ws.Range("G3").Value = 49
$ws.Range("H3").Value = 31
$ws.Range("I3").Value = 31
$ws.Range("K3").Value = 1045
$ws.Range("L3").Value = 374
$ws.Range("M3").Value = 671
$ws.Range("N3").Value = 671
$ws.Range("P3").Value = 128
$ws.Range("Q3").Value = 83
$ws.Range("R3").Value = -87
$ws.Range("S3").Value = 39
$ws.Range("T3").Value = 92
$ws.Range("U3").Value = -9
$ws.Range("V3").Value = 183
$ws.Range("W3").Value = 5.08
$ws.Range("X3").Value = 3.23
$ws.Range("Y3").Value = 4.62
$ws.Range("Z3").Value = 2.98
$ws.Range("AA3").Value = 55.76
$ws.Range("AB3").Value = 349.7
$ws.Range("AC3").Value = 118
$ws.Range("AD3").Value = 18.39
$ws.Range("AE3").Value = 2578
$ws.Range("AF3").Value = 0.84
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 26027444

# Row 4
$ws.Range("D4").Value = 1078
$ws.Range("E4").Value = 106
$ws.Range("F4").Value = 106
$ws.Range("G4").Value = 93
$ws.Range("H4").Value = 59
$ws.Range("I4").Value = 55
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 1656
$ws.Range("L4").Value = 861
$ws.Range("M4").Value = 794
$ws.Range("N4").Value = 762
$ws.Range("O4").Value = 32
$ws.Range("P4").Value = 131
$ws.Range("Q4").Value = -79
$ws.Range("R4").Value = -185
$ws.Range("S4").Value = 341
$ws.Range("T4").Value = 31
$ws.Range("U4").Value = -110
$ws.Range("V4").Value = 586
$ws.Range("W4").Value = 9.84
$ws.Range("X4").Value = 5.43
$ws.Range("Y4").Value = 7.74
$ws.Range("Z4").Value = 4.34
$ws.Range("AA4").Value = 108.43
$ws.Range("AB4").Value = 420.41
$ws.Range("AC4").Value = 212
$ws.Range("AD4").Value = 14.95
$ws.Range("AE4").Value = 2869
$ws.Range("AF4").Value = 1.11
$ws.Range("AG4").Value = 76
$ws.Range("AH4").Value = 2.39
$ws.Range("AI4").Value = 36.3
$ws.Range("AJ4").Value = 26552090

# Row 5
$ws.Range("D5").Value = 1120
$ws.Range("E5").Value = 84
$ws.Range("F5").Value = 84
$ws.Range("G5").Value = 63
$ws.Range("H5").Value = 37
$ws.Range("I5").Value = 43
$ws.Range("J5").Value = -6
$ws.Range("K5").Value = 1918
$ws.Range("L5").Value = 981
$ws.Range("M5").Value = 937
$ws.Range("N5").Value = 931
$ws.Range("O5").Value = 7
$ws.Range("P5").Value = 158
$ws.Range("Q5").Value = 88
$ws.Range("R5").Value = -101
$ws.Range("S5").Value = -42
$ws.Range("T5").Value = 58
$ws.Range("U5").Value = 30
$ws.Range("V5").Value = 631
$ws.Range("W5").Value = 7.5
$ws.Range("X5").Value = 3.33
$ws.Range("Y5").Value = 5.07
$ws.Range("Z5").Value = 2.08
$ws.Range("AA5").Value = 104.65
$ws.Range("AB5").Value = 449.72
$ws.Range("AC5").Value = 143
$ws.Range("AD5").Value = 22.31
$ws.Range("AE5").Value = 2971
$ws.Range("AF5").Value = 1.08
$ws.Range("AG5").Value = 49
$ws.Range("AH5").Value = 1.52
$ws.Range("AI5").Value = 35.45
$ws.Range("AJ5").Value = 32153118

# Row 6
$ws.Range("D6").Value = 1100
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = 10
$ws.Range("G6").Value = -102
$ws.Range("H6").Value = -110
$ws.Range("I6").Value = -100
$ws.Range("K6").Value = 3041
$ws.Range("L6").Value = 1441
$ws.Range("M6").Value = 1600
$ws.Range("N6").Value = 1595
$ws.Range("P6").Value = 271
$ws.Range("Q6").Value = -229
$ws.Range("R6").Value = -558
$ws.Range("S6").Value = 902
$ws.Range("T6").Value = 34
$ws.Range("U6").Value = -264
$ws.Range("V6").Value = 994
$ws.Range("W6").Value = 0.9399999999999999
$ws.Range("X6").Value = -10.04
$ws.Range("Y6").Value = -7.95
$ws.Range("Z6").Value = -4.45
$ws.Range("AA6").Value = 90.04000000000001
$ws.Range("AB6").Value = 508.64
$ws.Range("AC6").Value = -209
$ws.Range("AD6").Value = -58.32
$ws.Range("AE6").Value = 2963
$ws.Range("AF6").Value = 4.12
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 53826895

# Remove now-unused cells (columns no longer reported for these years)
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# Clear estimate-year rows (7-9): keep label columns, drop all reported figures
# Row 7
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()

